# Ambermoon Advanced workbook update: "Slightly empowered thief equip"
$wb = $excel.ActiveWorkbook

# --- Summary sheet: add a note about archive entry limits ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A1").Value = "Note: Never add more than 530 entries to an archive otherwise the Amiga version will fail."

# --- Items sheet: add "Changed Items" table describing balance changes ---
$items = $wb.Worksheets.Item("Items")

$items.Range("F1").Value = "Changed Items"
$items.Range("F1:H1").Merge()
$items.Range("F1").Style = $items.Range("A1").Style
$items.Range("G1:H1").Style = $items.Range("B1").Style

$items.Range("F2").Value = "Index"
$items.Range("G2").Value = "Name"
$items.Range("H2").Value = "Changes"

$items.Range("F3").Value = 249
$items.Range("G3").Value = "Shadow Belt"
$items.Range("H3").Value = "Number of charges (Blink) increased from 5 to 15"

$items.Range("F4").Value = 251
$items.Range("G4").Value = "Murder Blade"
$items.Range("H4").Value = "Atk increased from 10 to 11, MagicWeaponLevel from 0 to 1"

$items.Columns.Item(7).ColumnWidth = 16.140625

# The Items sheet becomes the active/selected sheet
$items.Range("H5").Select()
$items.Activate()

Write-Host "Applied thief equip balance changes"
